$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before the old "Estado" column (I),
# shifting Estado/Transaccion/Fecha/Cuenta from I:L to K:N.
$ws.Columns("I:J").Insert() | Out-Null

# New header row cells (style matches the rest of the blue header band, A1:H1).
$ws.Range("I1").Value = "metodo"
$ws.Range("J1").Value = "frecuencia"

# New data row cells (leading apostrophe keeps the quote-prefixed text style
# used by the other data cells in row 2, e.g. B2:H2). J2 is written first to
# match the shared-string insertion order of the original edit.
$ws.Range("J2").Value = "'e0Y e1M e0W o3D e0F"
$ws.Range("I2").Value = "'PAY"

# The old Estado/Transaccion/Fecha values (now in K2:M2) are removed entirely.
$ws.Range("K2:M2").ClearContents() | Out-Null

# Autofit the two new columns to their content, like Excel does after typing
# (best-fit to "metodo"/"PAY" and "frecuencia"/"e0Y e1M e0W o3D e0F").
$ws.Columns("I").ColumnWidth = 9.25
$ws.Columns("J").ColumnWidth = 18.6

# Move the selection, matching the end-state cursor position.
$ws.Range("L6").Select() | Out-Null
